$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 983.7273
$ws.Range("I12").Value = 546
$ws.Range("J12").Value = 1749.75
$ws.Range("K12").Value = 546
$ws.Range("L12").Value = 1749.75
$ws.Range("M12").Value = -376
$ws.Range("N12").Value = -2089.75
$ws.Range("H33").Value = 642.3
$ws.Range("J33").Value = 598
$ws.Range("L33").Value = 598
$ws.Range("N33").Value = -1056
$ws.Range("H52").Value = 6502.5884
$ws.Range("I52").Value = 1199.5
$ws.Range("J52").Value = 8134.3076
$ws.Range("K52").Value = 3598.5
$ws.Range("L52").Value = 24402.9228
$ws.Range("M52").Value = -3438.5
$ws.Range("N52").Value = -24722.9228
$ws.Range("H98").Value = 1308.9344
$ws.Range("I98").Value = 1325.0377
$ws.Range("K98").Value = 1325.0377
$ws.Range("M98").Value = 172.9622999999999
$ws.Range("H100").Value = 3500
$ws.Range("I100").Value = 4057.8572
$ws.Range("J100").Value = 2719
$ws.Range("K100").Value = 4057.8572
$ws.Range("L100").Value = 2719
$ws.Range("M100").Value = -3516.8572
$ws.Range("N100").Value = -3801
$ws.Range("H113").Value = 5849.476
$ws.Range("I113").Value = 5814.4
$ws.Range("J113").Value = 5881.364
$ws.Range("K113").Value = 5814.4
$ws.Range("L113").Value = 5881.364
$ws.Range("M113").Value = -2560.4
$ws.Range("N113").Value = -12389.364
$ws.Range("H115").Value = 923.6923
$ws.Range("I115").Value = 963.1667
$ws.Range("J115").Value = 450
$ws.Range("K115").Value = 2889.5001
$ws.Range("L115").Value = 1350
$ws.Range("M115").Value = -1322.5001
$ws.Range("N115").Value = -4484
$ws.Range("H116").Value = 5100.212
$ws.Range("I116").Value = 3987.2632
$ws.Range("K116").Value = 3987.2632
$ws.Range("M116").Value = -545.2631999999999
$ws.Range("H122").Value = 1308.9344
$ws.Range("I122").Value = 1325.0377
$ws.Range("K122").Value = 3975.1131
$ws.Range("M122").Value = -1525.1131

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 389.83334
$ws.Range("I5").Value = 353.7143
$ws.Range("J5").Value = 440.4
$ws.Range("K5").Value = 353.7143
$ws.Range("L5").Value = 440.4
$ws.Range("M5").Value = -241.7143
$ws.Range("N5").Value = -664.4
$ws.Range("H61").Value = 2326.08
$ws.Range("I61").Value = 2214.875
$ws.Range("K61").Value = 2214.875
$ws.Range("M61").Value = -2002.875
$ws.Range("H63").Value = 5721.933
$ws.Range("I63").Value = 2918.5715
$ws.Range("J63").Value = 8174.875
$ws.Range("K63").Value = 2918.5715
$ws.Range("L63").Value = 8174.875
$ws.Range("M63").Value = -2232.5715
$ws.Range("N63").Value = -9546.875
$ws.Range("H66").Value = 5721.933
$ws.Range("I66").Value = 2918.5715
$ws.Range("J66").Value = 8174.875
$ws.Range("K66").Value = 14592.8575
$ws.Range("L66").Value = 40874.375
$ws.Range("M66").Value = -11160.8575
$ws.Range("N66").Value = -47738.375
$ws.Range("H74").Value = 21017.66
$ws.Range("I74").Value = 1335.25
$ws.Range("K74").Value = 1335.25
$ws.Range("M74").Value = -461.25
$ws.Range("H77").Value = 21017.66
$ws.Range("I77").Value = 1335.25
$ws.Range("K77").Value = 6676.25
$ws.Range("M77").Value = -2308.25
$ws.Range("H136").Value = 2326.08
$ws.Range("I136").Value = 2214.875
$ws.Range("K136").Value = 6644.625
$ws.Range("M136").Value = -4094.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 389.83334
$ws.Range("I4").Value = 353.7143
$ws.Range("J4").Value = 440.4
$ws.Range("K4").Value = 353.7143
$ws.Range("L4").Value = 440.4
$ws.Range("M4").Value = -238.7143
$ws.Range("N4").Value = -670.4
$ws.Range("H86").Value = 10001500
$ws.Range("I86").Value = 20001300
$ws.Range("J86").Value = 1700.8
$ws.Range("K86").Value = 20001300
$ws.Range("L86").Value = 1700.8
$ws.Range("M86").Value = -20000177
$ws.Range("N86").Value = -3946.8
$ws.Range("H89").Value = 10001500
$ws.Range("I89").Value = 20001300
$ws.Range("J89").Value = 1700.8
$ws.Range("K89").Value = 100006500
$ws.Range("L89").Value = 8504
$ws.Range("M89").Value = -100000884
$ws.Range("N89").Value = -19736

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 690.25
$ws.Range("I7").Value = 509
$ws.Range("K7").Value = 509
$ws.Range("M7").Value = -396
$ws.Range("H31").Value = 22972.84
$ws.Range("I31").Value = 3010.3333
$ws.Range("K31").Value = 3010.3333
$ws.Range("M31").Value = -2715.3333
$ws.Range("H34").Value = 22972.84
$ws.Range("I34").Value = 3010.3333
$ws.Range("K34").Value = 3010.3333
$ws.Range("M34").Value = -2808.3333
$ws.Range("H99").Value = 3017.318
$ws.Range("I99").Value = 2622.7646
$ws.Range("J99").Value = 4358.8
$ws.Range("K99").Value = 2622.7646
$ws.Range("L99").Value = 4358.8
$ws.Range("M99").Value = -1124.7646
$ws.Range("N99").Value = -7354.8
$ws.Range("H126").Value = 3017.318
$ws.Range("I126").Value = 2622.7646
$ws.Range("J126").Value = 4358.8
$ws.Range("K126").Value = 7868.293799999999
$ws.Range("L126").Value = 13076.4
$ws.Range("M126").Value = -5398.293799999999
$ws.Range("N126").Value = -18016.4
$ws.Range("H134").Value = 2869.25
$ws.Range("I134").Value = 1980.9117
$ws.Range("K134").Value = 5942.7351
$ws.Range("M134").Value = -3407.7351

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 69500
$ws.Range("J37").Value = 69500
$ws.Range("L37").Value = 208500
$ws.Range("N37").Value = -208724
$ws.Range("H132").Value = 1726.5333
$ws.Range("I132").Value = 1513.8572
$ws.Range("J132").Value = 1912.625
$ws.Range("K132").Value = 13624.7148
$ws.Range("L132").Value = 17213.625
$ws.Range("M132").Value = -11094.7148
$ws.Range("N132").Value = -22273.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 9912.5
$ws.Range("J57").Value = 9912.5
$ws.Range("L57").Value = 9912.5
$ws.Range("N57").Value = -11552.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 7338.077
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 7338.077
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 7338.077
$ws.Range("M46").Value = $null
$ws.Range("N46").Value = -7714.077

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1318.6086
$ws.Range("I113").Value = 852.5714
$ws.Range("J113").Value = 1522.5
$ws.Range("K113").Value = 2557.7142
$ws.Range("L113").Value = 4567.5
$ws.Range("M113").Value = -387.7142000000003
$ws.Range("N113").Value = -8907.5
$ws.Range("H136").Value = 1902.6171
$ws.Range("I136").Value = 1467.0883
$ws.Range("K136").Value = 4401.2649
$ws.Range("M136").Value = -1851.2649
